$d = $word.ActiveDocument

# Locate the exact run of text that needs to be replaced:
#   setwd("C:/Users/Vasyl/Documents/GitHub/asSeq/pipeline_GTEx/v8/example/Muscle_Skeletal")
$q = [char]34
$oldText = "setwd(" + $q + "C:/Users/Vasyl/Documents/GitHub/asSeq/pipeline_GTEx/v8/example/Muscle_Skeletal" + $q + ")"

$finder = $d.Content
$finder.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$startPos = $finder.Start

# Build the full replacement text (both statements, no break yet) in a
# single Range.Text assignment so that sibling runs/line-breaks elsewhere
# in the paragraph are left untouched.
$part1 = "ts.dir =" + " " + "sprintf" + "(" + $q + "%s/Muscle_Skeletal" + $q + ", rt.dir)"
$part2 = "setwd" + "(ts.dir)"
$finder.Text = $part1 + $part2

# Re-derive each segment's range from $startPos and apply the matching
# syntax-highlighting character style to it (mirrors the rest of this
# "SourceCode" paragraph).
$p1 = $startPos
$p2 = $p1 + 8                      # "ts.dir ="
$d.Range($p1, $p2).Style = "NormalTok"

$p3 = $p2 + 1                      # " "
$d.Range($p2, $p3).Style = "StringTok"

$p4 = $p3 + 7                      # "sprintf"
$d.Range($p3, $p4).Style = "KeywordTok"

$p5 = $p4 + 1                      # "("
$d.Range($p4, $p5).Style = "NormalTok"

$p6 = $p5 + 20                     # "%s/Muscle_Skeletal" (incl. quotes)
$d.Range($p5, $p6).Style = "StringTok"

$p7 = $p6 + 9                      # ", rt.dir)"
$d.Range($p6, $p7).Style = "NormalTok"

$p8 = $p7 + 5                      # "setwd"
$d.Range($p7, $p8).Style = "KeywordTok"

$p9 = $p8 + 8                      # "(ts.dir)"
$d.Range($p8, $p9).Style = "NormalTok"

# Insert the manual line break between the two statements last, so it
# lands in its own run (matching how the document's other line breaks
# are structured) instead of being absorbed into a styled run.
$brk = $d.Range($p7, $p7)
$brk.InsertBreak(6)
